$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-17 05:18:59'
$ws.Range('M2').Value = '1.1 °C 4:56 TU'
$ws.Range('E3').Value = '2026-02-17 05:19:02'
$ws.Range('E4').Value = '2026-02-17 05:19:04'
$ws.Range('J4').Value = '1015.3 hPa'
$ws.Range('L4').Value = '12.6 km/h - 299º 4:48 TU'
$ws.Range('E5').Value = '2026-02-17 05:19:06'
$ws.Range('E6').Value = '2026-02-17 05:19:09'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '87%'
$ws.Range('J6').Value = '1015.0 hPa'
$ws.Range('O6').Value = '8.4 °C'
$ws.Range('E7').Value = '2026-02-17 05:19:11'
$ws.Range('H7').NumberFormat = '@'
$ws.Range('H7').Value = '54%'
$ws.Range('J7').Value = '1014.5 hPa'
$ws.Range('K7').Value = '-0.1 MJ/m2'
$ws.Range('N7').Value = '12.6 °C 4:41 TU'
$ws.Range('O7').Value = '14.2 °C'
$ws.Range('E8').Value = '2026-02-17 05:19:14'
$ws.Range('J8').Value = '1014.8 hPa'
$ws.Range('N8').Value = '8.9 °C 4:58 TU'
$ws.Range('O8').Value = '10.0 °C'
$ws.Range('E9').Value = '2026-02-17 05:19:17'
$ws.Range('H9').NumberFormat = '@'
$ws.Range('H9').Value = '45%'
$ws.Range('E10').Value = '2026-02-17 05:19:19'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '90%'
$ws.Range('O10').Value = '7.9 °C'
$ws.Range('E11').Value = '2026-02-17 05:19:22'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '37%'
$ws.Range('N11').Value = '2.4 °C 4:56 TU'
$ws.Range('O11').Value = '6.7 °C'
$ws.Range('E12').Value = '2026-02-17 05:19:24'
$ws.Range('H12').NumberFormat = '@'
$ws.Range('H12').Value = '46%'
$ws.Range('N12').Value = '11.6 °C 4:59 TU'
$ws.Range('O12').Value = '12.9 °C'
$ws.Range('E13').Value = '2026-02-17 05:19:27'
$ws.Range('J13').Value = '1016.6 hPa'
$ws.Range('O13').Value = '4.1 °C'
$ws.Range('E14').Value = '2026-02-17 05:19:29'
$ws.Range('N14').Value = '12.3 °C 4:59 TU'
$ws.Range('O14').Value = '13.4 °C'
$ws.Range('E15').Value = '2026-02-17 05:19:31'
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '46%'
$ws.Range('N15').Value = '11.2 °C 4:59 TU'
$ws.Range('O15').Value = '12.4 °C'
$ws.Range('E16').Value = '2026-02-17 05:19:34'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value = '45%'
$ws.Range('E17').Value = '2026-02-17 05:19:37'
$ws.Range('H17').NumberFormat = '@'
$ws.Range('H17').Value = '51%'
$ws.Range('E18').Value = '2026-02-17 05:19:39'
$ws.Range('J18').Value = '1015.3 hPa'
$ws.Range('E19').Value = '2026-02-17 05:19:42'
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '67%'
$ws.Range('E20').Value = '2026-02-17 05:19:43'
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '39%'
$ws.Range('E21').Value = '2026-02-17 05:19:44'
$ws.Range('J21').Value = '1015.5 hPa'
$ws.Range('K21').Value = '-0.1 MJ/m2'
$ws.Range('N21').Value = '4.8 °C 4:59 TU'
$ws.Range('O21').Value = '7.9 °C'
$ws.Range('E22').Value = '2026-02-17 05:19:45'
$ws.Range('E23').Value = '2026-02-17 05:19:46'
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '59%'
$ws.Range('L23').Value = '78.5 km/h - 328º 4:33 TU'
$ws.Range('E24').Value = '2026-02-17 05:19:47'
$ws.Range('J24').Value = '1017.4 hPa'
$ws.Range('E25').Value = '2026-02-17 05:19:48'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '43%'
$ws.Range('O25').Value = '-3.6 °C'
$ws.Range('E26').Value = '2026-02-17 05:19:51'
$ws.Range('E27').Value = '2026-02-17 05:19:54'
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '42%'
$ws.Range('N27').Value = '-3.5 °C 4:59 TU'
$ws.Range('O27').Value = '-2.6 °C'
$ws.Range('E28').Value = '2026-02-17 05:19:56'
$ws.Range('J28').Value = '1015.4 hPa'
$ws.Range('N28').Value = '3.7 °C 4:53 TU'
$ws.Range('O28').Value = '5.0 °C'
$ws.Range('E29').Value = '2026-02-17 05:19:59'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '60%'
$ws.Range('L29').Value = '24.5 km/h - 41º 4:40 TU'
$ws.Range('M29').Value = '13.0 °C 4:18 TU'
$ws.Range('O29').Value = '11.8 °C'
$ws.Range('E30').Value = '2026-02-17 05:20:02'
$ws.Range('J30').Value = '1014.4 hPa'
$ws.Range('N30').Value = '9.4 °C 4:41 TU'
$ws.Range('O30').Value = '12.0 °C'
$ws.Range('E31').Value = '2026-02-17 05:20:04'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '64%'
$ws.Range('J31').Value = '1015.1 hPa'
$ws.Range('N31').Value = '7.7 °C 4:36 TU'
$ws.Range('O31').Value = '9.9 °C'
$ws.Range('E32').Value = '2026-02-17 05:20:06'
$ws.Range('O32').Value = '6.4 °C'
$ws.Range('E33').Value = '2026-02-17 05:20:09'
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = '39%'
$ws.Range('J33').Value = '1015.7 hPa'
$ws.Range('E34').Value = '2026-02-17 05:20:12'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '45%'
$ws.Range('O34').Value = '-1.0 °C'
$ws.Range('E35').Value = '2026-02-17 05:20:14'
$ws.Range('E36').Value = '2026-02-17 05:20:16'
$ws.Range('H36').NumberFormat = '@'
$ws.Range('H36').Value = '48%'
$ws.Range('N36').Value = '10.6 °C 4:59 TU'
$ws.Range('O36').Value = '13.0 °C'
$ws.Range('E37').Value = '2026-02-17 05:20:19'
$ws.Range('H37').NumberFormat = '@'
$ws.Range('H37').Value = '42%'
$ws.Range('J37').Value = '1015.3 hPa'
$ws.Range('N37').Value = '4.9 °C 4:57 TU'
$ws.Range('O37').Value = '8.2 °C'
$ws.Range('E38').Value = '2026-02-17 05:20:22'
$ws.Range('O38').Value = '9.0 °C'
$ws.Range('E39').Value = '2026-02-17 05:20:24'
$ws.Range('H39').NumberFormat = '@'
$ws.Range('H39').Value = '55%'
$ws.Range('I39').Value = '1.0 mm'
$ws.Range('M39').Value = '-3.4 °C 4:57 TU'
$ws.Range('O39').Value = '-4.9 °C'
$ws.Range('E40').Value = '2026-02-17 05:20:27'
$ws.Range('O40').Value = '5.4 °C'
$ws.Range('E41').Value = '2026-02-17 05:20:29'
$ws.Range('J41').Value = '1015.2 hPa'
$ws.Range('N41').Value = '12.4 °C 4:59 TU'
$ws.Range('O41').Value = '15.1 °C'
$ws.Range('E42').Value = '2026-02-17 05:20:31'
$ws.Range('H42').NumberFormat = '@'
$ws.Range('H42').Value = '47%'
$ws.Range('O42').Value = '12.9 °C'
$ws.Range('E43').Value = '2026-02-17 05:20:34'
$ws.Range('N43').Value = '2.9 °C 4:32 TU'
$ws.Range('O43').Value = '4.5 °C'
$ws.Range('E44').Value = '2026-02-17 05:20:36'
$ws.Range('H44').NumberFormat = '@'
$ws.Range('H44').Value = '66%'
$ws.Range('E45').Value = '2026-02-17 05:20:39'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '51%'
$ws.Range('O45').Value = '4.6 °C'
$ws.Range('E46').Value = '2026-02-17 05:20:41'
$ws.Range('H46').NumberFormat = '@'
$ws.Range('H46').Value = '58%'
$ws.Range('J46').Value = '1017.6 hPa'
$ws.Range('N46').Value = '12.5 °C 4:34 TU'
$ws.Range('O46').Value = '13.5 °C'
